$wb = $excel.ActiveWorkbook

# --- Sheet "Livros": change D7 to a number and add a new row 8 ---
$wsLivros = $wb.Worksheets.Item("Livros")

$wsLivros.Range("D7").Value = 1234567890123

$wsLivros.Range("A8").Value = "aaadw"
$wsLivros.Range("B8").Value = "adawd"
$wsLivros.Range("C8").Value = "awdawd"
$wsLivros.Range("D8").NumberFormat = "@"
$wsLivros.Range("D8").Value = "1234567890123"
$wsLivros.Range("D8").ClearFormats()

# --- Sheet "Usuarios": add header row and a data row ---
$wsUsuarios = $wb.Worksheets.Item("Usuarios")

$wsUsuarios.Range("A1").Value = "nome"
$wsUsuarios.Range("B1").Value = "idade"
$wsUsuarios.Range("C1").Value = "cpf"
$wsUsuarios.Range("D1").Value = "email"
$wsUsuarios.Range("E1").Value = "telefone"
$wsUsuarios.Range("F1").Value = "endereco"

$wsLivros.Range("A1:D1").Copy()
$wsUsuarios.Range("A1:F1").PasteSpecial(-4122)

$wsUsuarios.Range("A2").Value = "abc"
$wsUsuarios.Range("B2").Value = 123
$wsUsuarios.Range("C2").NumberFormat = "@"
$wsUsuarios.Range("C2").Value = "12345678901"
$wsUsuarios.Range("C2").ClearFormats()
$wsUsuarios.Range("D2").Value = "abcde@gmail.com"
$wsUsuarios.Range("E2").NumberFormat = "@"
$wsUsuarios.Range("E2").Value = "1212341234"
$wsUsuarios.Range("E2").ClearFormats()
$wsUsuarios.Range("F2").Value = "rua1"
